$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the last existing data row (210) eight times, cascading the
# formatting (style) down into the new rows 211..218 before filling in
# their content.
$ws.Rows.Item(210).Copy()
$ws.Rows.Item(211).Insert()

$ws.Rows.Item(211).Copy()
$ws.Rows.Item(212).Insert()

$ws.Rows.Item(212).Copy()
$ws.Rows.Item(213).Insert()

$ws.Rows.Item(213).Copy()
$ws.Rows.Item(214).Insert()

$ws.Rows.Item(214).Copy()
$ws.Rows.Item(215).Insert()

$ws.Rows.Item(215).Copy()
$ws.Rows.Item(216).Insert()

$ws.Rows.Item(216).Copy()
$ws.Rows.Item(217).Insert()

$ws.Rows.Item(217).Copy()
$ws.Rows.Item(218).Insert()

# Fill in the new translation keys/values (column B = key, column C =
# translation). Row 218's content is entered before row 217's so the
# shared-string table is built up in the same order as the source file.
$ws.Range("B211").Value = "lab.atomizer.tooltip.create"
$ws.Range("C211").Value = "Přidat atomizér"

$ws.Range("B212").Value = "lab.atomizer.create.title"
$ws.Range("C212").Value = "Vytvořit atomizér"

$ws.Range("B213").Value = "lab.atomizer.create.subtitle"
$ws.Range("C213").Value = "Přidejte nový atomizér, který tak bude přístupný ostatních a v buildech."

$ws.Range("B214").Value = "lab.atomizer.name.label"
$ws.Range("C214").Value = "Jméno"

$ws.Range("B215").Value = "lab.atomizer.vendorId.label"
$ws.Range("C215").Value = "Výrobce"

$ws.Range("B216").Value = "lab.vendor.tooltip.create"
$ws.Range("C216").Value = "Založit výrobce"

$ws.Range("B218").Value = "lab.vendor.create.subtitle"
$ws.Range("C218").Value = "Výrobci jsou dostupní přes celou aplikaci, např. u atomizérů, modů, drátů a dalšího."

$ws.Range("B217").Value = "lab.vendor.create.title"
$ws.Range("C217").Value = "Založit výrobce"

# Match the author's final selection (cell the cursor ended up on).
$ws.Range("B211").Select()
